$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'39.849.06"
$ws.Range("E2").Value = "  -0.27%  "
$ws.Range("D3").Value = "'2.229.92"
$ws.Range("E3").Value = "  -4.38%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'298.55"
$ws.Range("E5").Value = "  -3.11%  "
$ws.Range("D6").Value = "'84.39"
$ws.Range("E6").Value = "  -0.44%  "
$ws.Range("D7").Value = "'0.514"
$ws.Range("E7").Value = "  -2.67%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").Value = "'0.468"
$ws.Range("E9").Value = "  -3.09%  "
$ws.Range("D10").Value = "'0.0781"
$ws.Range("E10").Value = "  -4.01%  "
$ws.Range("D11").Value = "'29.86"
$ws.Range("E11").Value = "  -0.65%  "
$ws.Range("D12").Value = "'46.71"
$ws.Range("E12").Value = "  -10.67%  "
$ws.Range("E13").Value = "  -2.14%  "
$ws.Range("D14").Value = "'2.575.26"
$ws.Range("E14").Value = "  -4.38%  "
$ws.Range("D15").Value = "'6.31"
$ws.Range("E15").Value = "  -1.70%  "
$ws.Range("D16").Value = "'14.15"
$ws.Range("E16").Value = "  -3.49%  "
$ws.Range("D17").Value = "'2.238.92"
$ws.Range("E17").Value = "  -4.33%  "
$ws.Range("D18").Value = "'0.719"
$ws.Range("E18").Value = "  -4.89%  "
$ws.Range("D19").Value = "'39.755.43"
$ws.Range("E19").Value = "  -0.41%  "
$ws.Range("D20").Value = "'0.0₃0880"
$ws.Range("E20").Value = "  -2.22%  "
$ws.Range("D21").Value = "'5.78"
$ws.Range("E21").Value = "  -5.27%  "
$ws.Range("B22").Value = "Litecoin"
$ws.Range("C22").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D22").Value = "'65.20"
$ws.Range("E22").Value = "  -3.80%  "
$ws.Range("B23").Value = "InternetComputer(DFINITY)"
$ws.Range("C23").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D23").Value = "'10.47"
$ws.Range("E23").Value = "  -1.11%  "
$ws.Range("D24").Value = "'234.74"
$ws.Range("E24").Value = "  -0.55%  "
$ws.Range("E25").Value = "  -0.08%  "
$ws.Range("E26").Value = "  -4.60%  "
$ws.Range("D27").Value = "'1.82"
$ws.Range("E27").Value = "  +0.61%  "
$ws.Range("D28").Value = "'22.81"
$ws.Range("E28").Value = "  -1.79%  "
$ws.Range("D29").Value = "'2.11"
$ws.Range("E29").Value = "  -0.41%  "
$ws.Range("D30").Value = "'9.21"
$ws.Range("E30").Value = "  -0.17%  "
$ws.Range("D31").Value = "'32.48"
$ws.Range("E31").Value = "  -4.41%  "
$ws.Range("D32").Value = "'149.68"
$ws.Range("E32").Value = "  -2.69%  "
$ws.Range("D34").Value = "'4.84"
$ws.Range("E34").Value = "  -4.89%  "
$ws.Range("E35").Value = "  -1.39%  "
$ws.Range("D36").Value = "'0.0702"
$ws.Range("E36").Value = "  -1.92%  "
$ws.Range("D37").Value = "'16.48"
$ws.Range("E37").Value = "  +6.33%  "
$ws.Range("E38").Value = "  -2.28%  "
$ws.Range("D39").Value = "'0.0978"
$ws.Range("E39").Value = "  -1.81%  "
$ws.Range("E40").Value = "  -2.05%  "
$ws.Range("E41").Value = "  -3.26%  "
$ws.Range("D42").Value = "'3.67"
$ws.Range("E42").Value = "  -3.20%  "
$ws.Range("D43").Value = "'1.933.33"
$ws.Range("E43").Value = "  -0.64%  "
$ws.Range("E44").Value = "  -2.86%  "
$ws.Range("D45").Value = "'0.0265"
$ws.Range("E45").Value = "  +1.35%  "
$ws.Range("D46").Value = "'9.25"
$ws.Range("E46").Value = "  +0.30%  "
$ws.Range("D47").Value = "'16.53"
$ws.Range("E47").Value = "  -5.66%  "
$ws.Range("E48").Value = "  -3.07%  "
$ws.Range("D49").Value = "'2.445.65"
$ws.Range("E49").Value = "  -4.15%  "
$ws.Range("D50").Value = "'71.25"
$ws.Range("E50").Value = "  +1.18%  "
$ws.Range("D51").Value = "'88.81"
$ws.Range("E51").Value = "  -3.95%  "
